$wb = $excel.ActiveWorkbook

# Sheet y1: update individual cell values
$ws1 = $wb.Worksheets.Item("y1")
$ws1.Range("B2").Value = 1
$ws1.Range("B4").Value = 0
$ws1.Range("A5").Value = 1
$ws1.Range("B6").Value = 0
$ws1.Range("A7").Value = 1
$ws1.Range("A8").Value = 1

# Sheet y2: update individual cell values
$ws2 = $wb.Worksheets.Item("y2")
$ws2.Range("A2").Value = 0
$ws2.Range("A5").Value = 0

# Sheet y3: update individual cell values
$ws3 = $wb.Worksheets.Item("y3")
$ws3.Range("A3").Value = 1
$ws3.Range("B3").Value = 0
$ws3.Range("A4").Value = 1
$ws3.Range("A6").Value = 1
$ws3.Range("A7").Value = 0
$ws3.Range("B8").Value = 0
